$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 87, shifting existing rows 87-185 down to 89-187
$ws.Rows("87:88").Insert()

# New row 87 values
$ws.Range("A87").Value = 9
$ws.Range("B87").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44539
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = 100112043
$ws.Range("G87").Value = "Pepino ensalada"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 106
$ws.Range("K87").Value = 6000
$ws.Range("L87").Value = 7000
$ws.Range("M87").Value = 6500
$ws.Range("N87").Value = "`$/caja 50 unidades"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 130
$ws.Range("Q87").Value = 50
$ws.Range("R87").Value = "Hortaliza"

# New row 88 values
$ws.Range("A88").Value = 9
$ws.Range("B88").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C88").Value = "Metropolitana"
$ws.Range("D88").Value = 44539
$ws.Range("E88").Value = 13
$ws.Range("F88").Value = 100112043
$ws.Range("G88").Value = "Pepino ensalada"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 160
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = 6500
$ws.Range("N88").Value = "`$/caja 50 unidades"
$ws.Range("O88").Value = "Región de O'Higgins"
$ws.Range("P88").Value = 130
$ws.Range("Q88").Value = 50
$ws.Range("R88").Value = "Hortaliza"
